# Config changed to low supervision: flip the supervision/control-point
# boolean switches (B5:B11 -> Sheet1!boolSupervision, boolControlPoint1-6)
# from TRUE to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5:B11").Value = $false

# Reflect the author's final selection/scroll position in the saved view
# (previously scrolled to B43 with B53 selected; now back at the top with
# B9 selected, matching the reviewed low-supervision switches).
$ws.Activate()
$ws.Range("B9").Select()
